$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 (Clojure) data values for columns C..I
$row19 = @{ "C19"=1; "D19"=2; "E19"=2; "F19"=3; "G19"=2; "H19"=3; "I19"=2 }
# Row 20 (Objective C) data values for columns C..I
$row20 = @{ "C20"=2; "D20"=2; "E20"=2; "F20"=1; "G20"=2; "H20"=2; "I20"=2 }

foreach ($addr in $row19.Keys) {
    $ws.Range($addr).Value = $row19[$addr]
}
foreach ($addr in $row20.Keys) {
    $ws.Range($addr).Value = $row20[$addr]
}

# Match the formatting of the row above (row 18) for the newly entered cells, including the
# empty "Completed?" cell in column J.
$ws.Range("C18:J18").Copy() | Out-Null
$ws.Range("C19:J19").PasteSpecial(-4122) | Out-Null
$ws.Range("C18:J18").Copy() | Out-Null
$ws.Range("C20:J20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update the selection to match the new active cell/selection.
$ws.Range("I22").Select() | Out-Null
